# Apply the edits described by the diff:
#   1. Fill in the CL1 (column C) reading-control scores for rows 15-24
#      (row 16 - CAICEDO POVEDA JEFFERSON PAUL - is left blank, matching
#      the diff). The G-column SUM formulas recalc automatically.
#   2. Select C8 on Hoja1 (its sheetView keeps tabSelected off, since it is
#      no longer the active tab after step 3).
#   3. Insert a new blank worksheet "Hoja2" right after "Hoja1" and leave
#      it as the active sheet/tab (workbookView activeTab + sheetView
#      tabSelected), matching the diff's new sheets entry + activeTab="1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 1. CL1 scores (column C), rows 15-24 ---
$ws.Range("C15").Value = 8
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 9
$ws.Range("C20").Value = 8
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 9
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 9

# --- 2. Selection on Hoja1 ---
$ws.Range("C8").Select()

# --- 3. New "Hoja2" worksheet, inserted after Hoja1, becomes active tab ---
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "Hoja2"
